$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "25.816.57"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.730.88"
$ws.Range("E3").Value = "  -1.87%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "229.16"
$ws.Range("E5").Value = "  -3.60%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.13%  "
Set-TextValue $ws.Range("D7") "0.5221"
$ws.Range("E7").Value = "  -0.79%  "
Set-TextValue $ws.Range("D8") "0.2752"
$ws.Range("E8").Value = "  +0.40%  "
Set-TextValue $ws.Range("D9") "39.27"
$ws.Range("E9").Value = "  -3.04%  "
Set-TextValue $ws.Range("D10") "0.06119"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "1.738.51"
$ws.Range("E11").Value = "  -1.45%  "
Set-TextValue $ws.Range("D12") "0.07055"
$ws.Range("E12").Value = "  +0.19%  "
Set-TextValue $ws.Range("D13") "14.94"
$ws.Range("E13").Value = "  -6.44%  "
Set-TextValue $ws.Range("D14") "0.6343"
$ws.Range("E14").Value = "  -3.29%  "
Set-TextValue $ws.Range("D15") "4.513"
$ws.Range("E15").Value = "  +0.12%  "
Set-TextValue $ws.Range("D16") "76.46"
$ws.Range("E16").Value = "  -2.64%  "
Set-TextValue $ws.Range("D17") "1.001"
$ws.Range("E17").Value = "  +0.22%  "
Set-TextValue $ws.Range("D18") "1.000"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "25.817.00"
$ws.Range("E19").Value = "  -0.96%  "
Set-TextValue $ws.Range("D20") "11.45"
$ws.Range("E20").Value = "  -2.31%  "
Set-TextValue $ws.Range("D21") "0.000006617"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").Value = "1.959.87"
$ws.Range("E22").Value = "  -1.75%  "
Set-TextValue $ws.Range("D23") "4.191"
$ws.Range("E23").Value = "  +2.30%  "
Set-TextValue $ws.Range("D24") "8.761"
$ws.Range("E24").Value = "  +3.98%  "
Set-TextValue $ws.Range("D25") "5.163"
$ws.Range("E25").Value = "  -0.70%  "
Set-TextValue $ws.Range("D26") "139.89"
$ws.Range("E26").Value = "  +1.42%  "
Set-TextValue $ws.Range("D27") "1.505"
$ws.Range("E27").Value = "  +1.37%  "
Set-TextValue $ws.Range("D28") "15.00"
$ws.Range("E28").Value = "  -1.58%  "
Set-TextValue $ws.Range("D29") "1.771"
$ws.Range("E29").Value = "  -3.81%  "
Set-TextValue $ws.Range("D30") "102.02"
$ws.Range("E30").Value = "  -1.09%  "
Set-TextValue $ws.Range("D31") "0.08266"
$ws.Range("E31").Value = "  -1.92%  "
Set-TextValue $ws.Range("D32") "3.703"
$ws.Range("E32").Value = "  -0.20%  "
Set-TextValue $ws.Range("D33") "3.495"
$ws.Range("E33").Value = "  +1.36%  "
Set-TextValue $ws.Range("D34") "0.04440"
$ws.Range("E34").Value = "  +0.17%  "
Set-TextValue $ws.Range("D35") "2.615"
$ws.Range("E35").Value = "  -1.30%  "
Set-TextValue $ws.Range("D36") "0.9680"
$ws.Range("E36").Value = "  -3.43%  "
Set-TextValue $ws.Range("D37") "0.6155"
$ws.Range("E37").Value = "  +0.87%  "
Set-TextValue $ws.Range("D38") "2.662"
$ws.Range("E38").Value = "  -3.17%  "
Set-TextValue $ws.Range("D39") "0.01563"
$ws.Range("E39").Value = "  -1.41%  "
Set-TextValue $ws.Range("D40") "1.000"
$ws.Range("E40").Value = "  -0.06%  "
Set-TextValue $ws.Range("D41") "1.898"
$ws.Range("E41").Value = "  -3.61%  "

$ws.Range("B42").Value = "PaxosStandard"
$ws.Range("C42").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
Set-TextValue $ws.Range("D42") "1.001"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D43") "99.42"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D44") "0.3807"
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D45") "5.009"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D46") "0.7179"
$ws.Range("E46").Value = "  -4.83%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.05331"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D48") "0.1115"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D49") "6.144"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D50") "53.09"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D51") "29.86"
$ws.Range("E51").Value = "  -1.27%  "
